$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-10
# from 2023-09-06 (45175) to 2023-09-14 (45183)
$ws.Range("C2:C10").Value = 45183
